# Slide 1, "TextBox 13" (shape 7 in z-order), paragraph "COLLEGE: ..." --
# the run that used to read " arts and science College /" must become
# three runs:
#   " arts and science College, "  (unchanged rPr)
#   "Bommayapalayam"                (unchanged rPr, flagged err="1" by PPT)
#   " /"                            (unchanged rPr)
# i.e. a comma + the word "Bommayapalayam" is inserted between
# "College" and the trailing " /".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(7)
$tr = $shp.TextFrame.TextRange

$oldFragment = " arts and science College /"
$word = "Bommayapalayam"
$leftPart = " arts and science College, "
$rightPart = " /"
$newFragment = $leftPart + $word + $rightPart

# --- locate the original run's text in the full story ---------------------
$full = $tr.Text
$start0 = $full.IndexOf($oldFragment)
if ($start0 -lt 0) {
    throw "Could not find target fragment in TextBox 13"
}

# The single space that separates "College" and "/" inside the run -- an
# interior character (not the run's first or last char), so replacing it
# forces the COM layer to split the run into three pieces instead of just
# rewriting the whole run's <a:t> in place.
$collegeHead = " arts and science College"
$spacePos = $start0 + $collegeHead.Length + 1   # 1-indexed position of that space
$spaceRange = $tr.Characters($spacePos, 1)
Write-Host ("step0 spaceRange=[" + $spaceRange.Text + "]")
$spaceRange.Text = ", " + $word + " "
Write-Host ("step1=[" + $tr.Text + "]")

# --- shrink the new middle run down to just the bare word ------------------
$full = $tr.Text
$midNeedle = ", " + $word + " "
$midStart0 = $full.IndexOf($midNeedle)
$wordPos = $midStart0 + 2 + 1                   # skip the leading ", "
$wordRange = $tr.Characters($wordPos, $word.Length)
Write-Host ("step1b wordRange=[" + $wordRange.Text + "]")
$wordRange.Text = $word
Write-Host ("step2=[" + $tr.Text + "]")

# --- merge the leading ", " back onto the "College" run ---------------------
$full = $tr.Text
$leftStart0 = $full.IndexOf($leftPart)
$leftRange = $tr.Characters($leftStart0 + 1, $leftPart.Length)
Write-Host ("step2b leftRange=[" + $leftRange.Text + "]")
$leftRange.Text = $leftPart
Write-Host ("step3=[" + $tr.Text + "]")

# --- merge the trailing " " back onto the "/" run ---------------------------
$full = $tr.Text
$wIdx0 = $full.IndexOf($word)
$tailPos = $wIdx0 + $word.Length + 1
$tailRange = $tr.Characters($tailPos, $rightPart.Length)
Write-Host ("step3b tailRange=[" + $tailRange.Text + "]")
$tailRange.Text = $rightPart
Write-Host ("step4=[" + $tr.Text + "]")

# --- sanity check ------------------------------------------------------------
$final = $tr.Text
if ($final.IndexOf($newFragment) -lt 0) {
    throw "Post-condition failed: expected fragment not present"
}
Write-Host ("OK: " + $newFragment)
